$d = $word.ActiveDocument

# Table 1 (Planning Poker results) - update "Puntos" column values
$t1 = $d.Tables.Item(1)
$t1.Cell(2,2).Range.Text = "4"   # Registro cliente/artista: 2 -> 4
$t1.Cell(3,2).Range.Text = "6"   # Inicio de sesión cliente/artista: 2 -> 6
$t1.Cell(4,2).Range.Text = "3"   # Filtros del catálogo: 2 -> 3
$t1.Cell(5,2).Range.Text = "7"   # Agregar camiseta al carrito: 6 -> 7
$t1.Cell(6,2).Range.Text = "7"   # Personalizar camiseta: 6 -> 7
$t1.Cell(7,2).Range.Text = "7"   # Eliminar Camisetas del carrito: 4 -> 7
$t1.Cell(8,2).Range.Text = "5"   # Ver productos del carrito: 4 -> 5
$t1.Cell(9,2).Range.Text = "9"   # Pagar camisetas: 8 -> 9

# Table 2 (Planning Poker results) - update "Puntos" column values
$t2 = $d.Tables.Item(2)
$t2.Cell(2,2).Range.Text = "6"   # Crear catálogo: 5 -> 6
$t2.Cell(3,2).Range.Text = "7"   # Añadir diseño: 4 -> 7
$t2.Cell(4,2).Range.Text = "6"   # Eliminar diseño: 3 -> 6
$t2.Cell(5,2).Range.Text = "4"   # Habilitar producto para venta: 3 -> 4
$t2.Cell(6,2).Range.Text = "7"   # Ver estadísticas de ventas: 5 -> 7

# Table 3 (Planning Poker results) - update "Puntos" column values
$t3 = $d.Tables.Item(3)
$t3.Cell(2,2).Range.Text = "6"   # Ver rating de productos: 4 -> 6
$t3.Cell(3,2).Range.Text = "6"   # Consultar y actualizar tarifas: 5 -> 6
